$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.505.95'
$ws.Range('E2').Value = '  +3.00%  '
$ws.Range('D3').Value = '3.366.21'
$ws.Range('E3').Value = '  +4.39%  '
$ws.Range('D5').Value = '''192.13'
$ws.Range('E5').Value = '  +4.70%  '
$ws.Range('D6').Value = '''593.21'
$ws.Range('E6').Value = '  +2.59%  '
$ws.Range('D7').Value = '''0.610'
$ws.Range('E7').Value = '  +1.03%  '
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('E9').Value = '  +3.29%  '
$ws.Range('D10').Value = '''6.77'
$ws.Range('E10').Value = '  +3.14%  '
$ws.Range('E11').Value = '  +2.57%  '
$ws.Range('D12').Value = '3.956.68'
$ws.Range('E12').Value = '  +4.67%  '
$ws.Range('E13').Value = '  +1.34%  '
$ws.Range('D14').Value = '''28.69'
$ws.Range('E14').Value = '  +3.33%  '
$ws.Range('D15').Value = '69.521.87'
$ws.Range('E15').Value = '  +2.94%  '
$ws.Range('E16').Value = '  +2.09%  '
$ws.Range('D17').Value = '3.354.50'
$ws.Range('E17').Value = '  +4.40%  '
$ws.Range('D18').Value = '''450.90'
$ws.Range('E18').Value = '  +14.06%  '
$ws.Range('E19').Value = '  +1.82%  '
$ws.Range('E20').Value = '  +3.18%  '
$ws.Range('E21').Value = '  +3.72%  '
$ws.Range('D22').Value = '''74.48'
$ws.Range('E22').Value = '  +4.68%  '
$ws.Range('D23').Value = '''0.999'
$ws.Range('E23').Value = '  -0.13%  '
$ws.Range('D24').Value = '3.517.84'
$ws.Range('E24').Value = '  +4.62%  '
$ws.Range('D25').Value = '''0.520'
$ws.Range('E25').Value = '  +1.10%  '
$ws.Range('E26').Value = '  +3.84%  '
$ws.Range('E27').Value = '  +4.61%  '
$ws.Range('E28').Value = '  +0.23%  '
$ws.Range('E29').Value = '  +0.09%  '
$ws.Range('E30').Value = '  +2.56%  '
$ws.Range('D31').Value = '''23.29'
$ws.Range('E31').Value = '  +2.93%  '
$ws.Range('E32').Value = '  +0.88%  '
$ws.Range('E33').Value = '  +3.99%  '
$ws.Range('E34').Value = '  +1.45%  '
$ws.Range('E35').Value = '  -0.01%  '
$ws.Range('E36').Value = '  +3.66%  '
$ws.Range('D37').Value = '''165.15'
$ws.Range('E37').Value = '  +2.66%  '
$ws.Range('E38').Value = '  +3.71%  '
$ws.Range('D39').Value = '''27.22'
$ws.Range('E39').Value = '  +3.72%  '
$ws.Range('E40').Value = '  +2.18%  '
$ws.Range('E41').Value = '  +1.31%  '
$ws.Range('D42').Value = '''6.56'
$ws.Range('E42').Value = '  +1.17%  '
$ws.Range('D43').Value = '2.740.75'
$ws.Range('E43').Value = '  +5.77%  '
$ws.Range('E44').Value = '  +3.12%  '
$ws.Range('D45').Value = '''25.72'
$ws.Range('E45').Value = '  +4.99%  '
$ws.Range('D46').Value = '''0.0690'
$ws.Range('E46').Value = '  +0.98%  '
$ws.Range('D47').Value = '''343.28'
$ws.Range('E47').Value = '  +2.89%  '
$ws.Range('D48').Value = '''40.78'
$ws.Range('E48').Value = '  +0.65%  '
$ws.Range('E49').Value = '  +3.14%  '
$ws.Range('B50').Value = 'ONDO'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D50').Value = '''1.04'
$ws.Range('E50').Value = '  +7.92%  '
$ws.Range('B51').Value = 'Arweave'
$ws.Range('C51').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D51').Value = '''33.05'
$ws.Range('E51').Value = '  +8.01%  '
